$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture a Style reference that already carries the built-in "Hyperlink" cell
# format (same xf index used by F2/G2) BEFORE any edits, so we can re-apply it
# after re-creating the hyperlinks below (Hyperlinks.Add otherwise mints a
# fresh, redundant style entry).
$hyperlinkStyle = $ws.Range("F2").Style

# --- Row 2 form data -------------------------------------------------
$ws.Range("A2").Value = "Sandeep"
$ws.Range("B2").Value = "sandeep@gmail.com"
$ws.Range("C2").Value = 7817008251
$ws.Range("D2").Value = "Vadodara"
$ws.Range("F2").Value = "Sandeep@123"
$ws.Range("G2").Value = "Sandeep@123"

# --- Hyperlinks --------------------------------------------------------
# This runtime has no in-place "update target" call, so drop the whole
# collection and rebuild it (in the same order/cells) pointing at the new
# mailto targets.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:sandeep@gmail.com")
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:Sandeep@123")
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:Sandeep@123")

# Restore the original Hyperlink cell style on the three linked cells.
$ws.Range("B2").Style = $hyperlinkStyle
$ws.Range("F2").Style = $hyperlinkStyle
$ws.Range("G2").Style = $hyperlinkStyle

# --- Selection -----------------------------------------------------
$ws.Range("C2").Select() | Out-Null
